$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=7; I='ba'; J='Appreciation'}
    @{Row=8; I='sd'; J='Statement-non-opinion'}
    @{Row=16; I='sd'; J='Statement-non-opinion'}
    @{Row=25; I='sd'; J='Statement-non-opinion'}
    @{Row=35; I='sd'; J='Statement-non-opinion'}
    @{Row=48; I='sd'; J='Statement-non-opinion'}
    @{Row=51; I='%'; J='Uninterpretable'}
    @{Row=52; I='sd'; J='Statement-non-opinion'}
    @{Row=53; I='sd'; J='Statement-non-opinion'}
    @{Row=56; I='sd'; J='Statement-non-opinion'}
    @{Row=57; I='sv'; J='Statement-opinion'}
    @{Row=68; I='sv'; J='Statement-opinion'}
    @{Row=86; I='sd'; J='Statement-non-opinion'}
    @{Row=91; I='sd'; J='Statement-non-opinion'}
    @{Row=92; I='sd'; J='Statement-non-opinion'}
    @{Row=113; I='%'; J='Uninterpretable'}
    @{Row=128; I='%'; J='Uninterpretable'}
    @{Row=129; I='aa'; J='Agree/Accept'}
    @{Row=145; I='sd'; J='Statement-non-opinion'}
    @{Row=148; I='sd'; J='Statement-non-opinion'}
    @{Row=155; I='sd'; J='Statement-non-opinion'}
    @{Row=158; I='qy'; J='Yes-No-Question'}
    @{Row=169; I='qy'; J='Yes-No-Question'}
    @{Row=183; I='aa'; J='Agree/Accept'}
    @{Row=184; I='aa'; J='Agree/Accept'}
    @{Row=195; I='sd'; J='Statement-non-opinion'}
    @{Row=213; I='sv'; J='Statement-opinion'}
    @{Row=222; I='sv'; J='Statement-opinion'}
    @{Row=229; I='sd'; J='Statement-non-opinion'}
    @{Row=235; I='sv'; J='Statement-opinion'}
    @{Row=238; I='sd'; J='Statement-non-opinion'}
    @{Row=245; I='sd'; J='Statement-non-opinion'}
    @{Row=256; I='ba'; J='Appreciation'}
    @{Row=259; I='sd'; J='Statement-non-opinion'}
    @{Row=267; I='sd'; J='Statement-non-opinion'}
    @{Row=276; I='sv'; J='Statement-opinion'}
    @{Row=283; I='ba'; J='Appreciation'}
    @{Row=286; I='sd'; J='Statement-non-opinion'}
    @{Row=295; I='sd'; J='Statement-non-opinion'}
    @{Row=300; I='aa'; J='Agree/Accept'}
    @{Row=305; I='sv'; J='Statement-opinion'}
    @{Row=307; I='sv'; J='Statement-opinion'}
    @{Row=319; I='sv'; J='Statement-opinion'}
    @{Row=326; I='aa'; J='Agree/Accept'}
    @{Row=341; I='sv'; J='Statement-opinion'}
    @{Row=347; I='sv'; J='Statement-opinion'}
    @{Row=356; I='sv'; J='Statement-opinion'}
    @{Row=373; I='sd'; J='Statement-non-opinion'}
    @{Row=397; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=429; I='ba'; J='Appreciation'}
    @{Row=430; I='sd'; J='Statement-non-opinion'}
    @{Row=432; I='ba'; J='Appreciation'}
    @{Row=437; I='sv'; J='Statement-opinion'}
    @{Row=438; I='sd'; J='Statement-non-opinion'}
    @{Row=451; I='sd'; J='Statement-non-opinion'}
    @{Row=453; I='%'; J='Uninterpretable'}
    @{Row=461; I='sd'; J='Statement-non-opinion'}
    @{Row=462; I='%'; J='Uninterpretable'}
    @{Row=486; I='aa'; J='Agree/Accept'}
    @{Row=500; I='%'; J='Uninterpretable'}
    @{Row=506; I='%'; J='Uninterpretable'}
    @{Row=507; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=510; I='%'; J='Uninterpretable'}
    @{Row=511; I='sv'; J='Statement-opinion'}
    @{Row=516; I='sd'; J='Statement-non-opinion'}
    @{Row=517; I='sv'; J='Statement-opinion'}
    @{Row=522; I='aa'; J='Agree/Accept'}
    @{Row=530; I='sd'; J='Statement-non-opinion'}
    @{Row=543; I='sd'; J='Statement-non-opinion'}
    @{Row=546; I='aa'; J='Agree/Accept'}
    @{Row=549; I='ba'; J='Appreciation'}
    @{Row=551; I='sd'; J='Statement-non-opinion'}
    @{Row=555; I='aa'; J='Agree/Accept'}
    @{Row=556; I='aa'; J='Agree/Accept'}
    @{Row=558; I='aa'; J='Agree/Accept'}
    @{Row=655; I='aa'; J='Agree/Accept'}
    @{Row=656; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=658; I='b'; J='Acknowledge (Backchannel)'}
    @{Row=680; I='sd'; J='Statement-non-opinion'}
    @{Row=706; I='%'; J='Uninterpretable'}
    @{Row=712; I='%'; J='Uninterpretable'}
    @{Row=716; I='sd'; J='Statement-non-opinion'}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

"Updated $($updates.Count) rows"
